# Auto-generated script to apply cryptos.xlsx data refresh diff
# (crypto price/volume table update, commit: "Updated cryptos list on Thu Jun 20 18:41:20 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price values look like plain decimal numbers (e.g. "593.85").
# Excel would silently coerce those into floating point numbers on assignment,
# which would both lose the exact textual formatting and diverge from the
# original text-valued cells. Force those specific cells to Text format first
# so the assigned strings are preserved verbatim.
$numericTextCells = @("D5", "D6", "D9", "D11", "D12", "D14", "D15", "D17", "D19", "D20", "D21", "D22", "D25", "D30", "D32", "D34", "D37", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all cell value updates from the diff
$ws.Range("D2").Value = "65.054.76"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.528.68"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "593.85"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").Value = "134.66"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("D7").Value = "3.524.90"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").Value = "7.13"
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "4.132.00"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "27.70"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "0.0000182"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.533.56"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.117"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "65.095.84"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "10.09"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "14.44"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "5.71"
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("D22").Value = "393.34"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").Value = "3.674.98"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").Value = "74.65"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -3.95%  "
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  +10.81%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").Value = "8.39"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").Value = "3.533.87"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("D34").Value = "24.22"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "5.30"
$ws.Range("E37").Value = "  +5.70%  "
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").Value = "1.57"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("D40").Value = "168.82"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").Value = "0.0816"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("D42").Value = "0.826"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("E43").Value = "  +4.03%  "
$ws.Range("D44").Value = "25.92"
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("D45").Value = "42.92"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "4.43"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "1.65"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "6.93"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").Value = "2.418.57"
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("E51").Value = "  +5.84%  "
